$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Rspo3-Lgr6, FAPs -> FAPs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.822099333333334
$ws.Range("H2").Value = 8.466298
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.2280253333333333
$ws.Range("N2").Value = 0.684076
$ws.Range("O2").Value = 0.5347665151402629
$ws.Range("P2").Value = 0.5347665151402629
$ws.Range("Q2").Value = 0.6435101411831111
$ws.Range("R2").Value = 5.791591270648
$ws.Range("S2").Value = 0.5347665151402629
$ws.Range("T2").Value = 0.5347665151402629

# Row 3 (Rspo3-Lgr6, FAPs -> sCs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.822099333333334
$ws.Range("H3").Value = 8.466298
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.1983763333333333
$ws.Range("N3").Value = 0.595129
$ws.Range("O3").Value = 0.4652334848597371
$ws.Range("P3").Value = 0.4652334848597371
$ws.Range("Q3").Value = 0.5598377180491112
$ws.Range("R3").Value = 5.038539462442
$ws.Range("S3").Value = 0.4652334848597371
$ws.Range("T3").Value = 0.4652334848597371
